$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (FAPs -> Rspo4 -> Lgr5 -> FAPs) picks up the recomputed values that
# used to live on row 3, but with updated (new-TPM) specificity figures.
$ws.Range("D2").Value = "FAPs"
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.3961209999999999
$ws.Range("N2").Value = 1.188363
$ws.Range("O2").Value = 0.5646784620538419
$ws.Range("P2").Value = 0.5646784620538419
$ws.Range("Q2").Value = 0.01881614362099999
$ws.Range("R2").Value = 0.169345292589
$ws.Range("S2").Value = 0.5646784620538419
$ws.Range("T2").Value = 0.5646784620538419

# Row 3 (FAPs -> Rspo4 -> Lgr5 -> MuSCs) picks up the recomputed values that
# used to live on row 4.
$ws.Range("D3").Value = "MuSCs"
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.3053773333333333
$ws.Range("N3").Value = 0.9161319999999999
$ws.Range("O3").Value = 0.4353215379461581
$ws.Range("P3").Value = 0.4353215379461581
$ws.Range("Q3").Value = 0.01450572871066667
$ws.Range("R3").Value = 0.130551558396
$ws.Range("S3").Value = 0.4353215379461581
$ws.Range("T3").Value = 0.4353215379461581

# Old row 4 (the MuSCs row) is now fully absorbed into row 3 above, so the
# trailing row is removed (dimension shrinks from A1:T4 to A1:T3, and the
# now-unused "ECs" shared string disappears on save).
$ws.Rows.Item(4).Delete()
